$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the bold/custom formatting that was applied to row 2 (axi_dma header row)
$ws.Rows.Item(2).ClearFormats()

# Swap the descriptions for the nfft (H5) and nfft_scaled (H8) rows
$h5 = $ws.Range("H5").Value2
$h8 = $ws.Range("H8").Value2
$ws.Range("H5").Value = $h8
$ws.Range("H8").Value = $h5

# Move the active selection to reflect where the user last clicked
$ws.Range("F33").Select()
